$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 157.8057217802531
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 407.9498229744228

$ws.Range("B3").Value = 0.001754667048134761
$ws.Range("C3").Value = 1766.335244827366
$ws.Range("D3").Value = 10137753.70137369
$ws.Range("E3").Value = 71517.89157740913
$ws.Range("G3").Value = 10211037.9299506
